# Insert a new daily price record for "Alcachofa" (Vega Central Mapocho de
# Santiago) at row 311. Excel's native row-insert semantics push the
# existing row 311 (and everything below it) down by one row, which is
# exactly the shift shown in the diff: the new record occupies row 311
# while the old rows 311..383 become rows 312..384.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 311, shifting rows 311:383 -> 312:384
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A311").Value = 9
$ws.Range("B311").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C311").Value = 'Metropolitana'
$ws.Range("D311").Value = 44694
$ws.Range("E311").Value = 13
$ws.Range("F311").Value = 100112013
$ws.Range("G311").Value = 'Alcachofa'
$ws.Range("H311").Value = 'Española'
$ws.Range("I311").Value = 'Primera'
$ws.Range("J311").Value = 52
$ws.Range("K311").Value = 24000
$ws.Range("L311").Value = 25000
$ws.Range("M311").Value = 24500
$ws.Range("N311").Value = '$/caja 30 unidades'
$ws.Range("O311").Value = 'Provincia del Elquí'
$ws.Range("P311").Value = 817
$ws.Range("Q311").Value = 30
$ws.Range("R311").Value = 'Hortaliza'
